$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 16,20

$arr[0,0] = "ECs"
$arr[0,1] = "Lpl"
$arr[0,2] = "Sdc1"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 46.48074800000001
$arr[0,7] = 139.442244
$arr[0,8] = 0.1473944418036112
$arr[0,9] = 0.1473944418036112
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.9159313333333333
$arr[0,13] = 2.747794
$arr[0,14] = 0.03641350786393945
$arr[0,15] = 0.03641350786393944
$arr[0,16] = 42.57317348997067
$arr[0,17] = 383.158561409736
$arr[0,18] = 0.005367148665716761
$arr[0,19] = 0.005367148665716759

$arr[1,0] = "ECs"
$arr[1,1] = "Lpl"
$arr[1,2] = "Sdc1"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 46.48074800000001
$arr[1,7] = 139.442244
$arr[1,8] = 0.1473944418036112
$arr[1,9] = 0.1473944418036112
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 10.56834633333333
$arr[1,13] = 31.705039
$arr[1,14] = 0.4201521973455823
$arr[1,15] = 0.4201521973455822
$arr[1,16] = 491.2246426963907
$arr[1,17] = 4421.021784267517
$arr[1,18] = 0.06192809860031278
$arr[1,19] = 0.06192809860031277

$arr[2,0] = "ECs"
$arr[2,1] = "Lpl"
$arr[2,2] = "Sdc1"
$arr[2,3] = "M2"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 46.48074800000001
$arr[2,7] = 139.442244
$arr[2,8] = 0.1473944418036112
$arr[2,9] = 0.1473944418036112
$arr[2,10] = 2
$arr[2,11] = 0.6666666666666666
$arr[2,12] = 1.220967333333333
$arr[2,13] = 3.662902
$arr[2,14] = 0.04854043308262539
$arr[2,15] = 0.04854043308262537
$arr[2,16] = 56.75147493689867
$arr[2,17] = 510.7632744320881
$arr[2,18] = 0.00715459003911911
$arr[2,19] = 0.007154590039119108

$arr[3,0] = "ECs"
$arr[3,1] = "Lpl"
$arr[3,2] = "Sdc1"
$arr[3,3] = "sCs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 46.48074800000001
$arr[3,7] = 139.442244
$arr[3,8] = 0.1473944418036112
$arr[3,9] = 0.1473944418036112
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 12.44836933333333
$arr[3,13] = 37.345108
$arr[3,14] = 0.494893861707853
$arr[3,15] = 0.4948938617078529
$arr[3,16] = 578.6095179935947
$arr[3,17] = 5207.485661942353
$arr[3,18] = 0.07294460449846253
$arr[3,19] = 0.07294460449846252

$arr[4,0] = "FAPs"
$arr[4,1] = "Lpl"
$arr[4,2] = "Sdc1"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 192.830597
$arr[4,7] = 578.4917909999999
$arr[4,8] = 0.6114823756165045
$arr[4,9] = 0.6114823756165044
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 0.9159313333333333
$arr[4,13] = 2.747794
$arr[4,14] = 0.03641350786393945
$arr[4,15] = 0.03641350786393944
$arr[4,16] = 176.6195858176727
$arr[4,17] = 1589.576272359054
$arr[4,18] = 0.02226621829317197
$arr[4,19] = 0.02226621829317195

$arr[5,0] = "FAPs"
$arr[5,1] = "Lpl"
$arr[5,2] = "Sdc1"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 192.830597
$arr[5,7] = 578.4917909999999
$arr[5,8] = 0.6114823756165045
$arr[5,9] = 0.6114823756165044
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 10.56834633333333
$arr[5,13] = 31.705039
$arr[5,14] = 0.4201521973455823
$arr[5,15] = 0.4201521973455822
$arr[5,16] = 2037.900532759427
$arr[5,17] = 18341.10479483485
$arr[5,18] = 0.2569156637533711
$arr[5,19] = 0.256915663753371

$arr[6,0] = "FAPs"
$arr[6,1] = "Lpl"
$arr[6,2] = "Sdc1"
$arr[6,3] = "M2"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 192.830597
$arr[6,7] = 578.4917909999999
$arr[6,8] = 0.6114823756165045
$arr[6,9] = 0.6114823756165044
$arr[6,10] = 2
$arr[6,11] = 0.6666666666666666
$arr[6,12] = 1.220967333333333
$arr[6,13] = 3.662902
$arr[6,14] = 0.04854043308262539
$arr[6,15] = 0.04854043308262537
$arr[6,16] = 235.4398598041646
$arr[6,17] = 2118.958738237482
$arr[6,18] = 0.02968161933481774
$arr[6,19] = 0.02968161933481772

$arr[7,0] = "FAPs"
$arr[7,1] = "Lpl"
$arr[7,2] = "Sdc1"
$arr[7,3] = "sCs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 192.830597
$arr[7,7] = 578.4917909999999
$arr[7,8] = 0.6114823756165045
$arr[7,9] = 0.6114823756165044
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 12.44836933333333
$arr[7,13] = 37.345108
$arr[7,14] = 0.494893861707853
$arr[7,15] = 0.4948938617078529
$arr[7,16] = 2400.426490223158
$arr[7,17] = 21603.83841200843
$arr[7,18] = 0.3026188742351438
$arr[7,19] = 0.3026188742351437

$arr[8,0] = "M2"
$arr[8,1] = "Lpl"
$arr[8,2] = "Sdc1"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 60.030993
$arr[8,7] = 180.092979
$arr[8,8] = 0.1903634318482028
$arr[8,9] = 0.1903634318482028
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 0.9159313333333333
$arr[8,13] = 2.747794
$arr[8,14] = 0.03641350786393945
$arr[8,15] = 0.03641350786393944
$arr[8,16] = 54.984267459814
$arr[8,17] = 494.858407138326
$arr[8,18] = 0.006931800322611035
$arr[8,19] = 0.006931800322611032

$arr[9,0] = "M2"
$arr[9,1] = "Lpl"
$arr[9,2] = "Sdc1"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 60.030993
$arr[9,7] = 180.092979
$arr[9,8] = 0.1903634318482028
$arr[9,9] = 0.1903634318482028
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 10.56834633333333
$arr[9,13] = 31.705039
$arr[9,14] = 0.4201521973455823
$arr[9,15] = 0.4201521973455822
$arr[9,16] = 634.428324757909
$arr[9,17] = 5709.854922821181
$arr[9,18] = 0.07998161418526842
$arr[9,19] = 0.07998161418526838

$arr[10,0] = "M2"
$arr[10,1] = "Lpl"
$arr[10,2] = "Sdc1"
$arr[10,3] = "M2"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 60.030993
$arr[10,7] = 180.092979
$arr[10,8] = 0.1903634318482028
$arr[10,9] = 0.1903634318482028
$arr[10,10] = 2
$arr[10,11] = 0.6666666666666666
$arr[10,12] = 1.220967333333333
$arr[10,13] = 3.662902
$arr[10,14] = 0.04854043308262539
$arr[10,15] = 0.04854043308262537
$arr[10,16] = 73.295881440562
$arr[10,17] = 659.6629329650581
$arr[10,18] = 0.009240323425006606
$arr[10,19] = 0.009240323425006603

$arr[11,0] = "M2"
$arr[11,1] = "Lpl"
$arr[11,2] = "Sdc1"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 60.030993
$arr[11,7] = 180.092979
$arr[11,8] = 0.1903634318482028
$arr[11,9] = 0.1903634318482028
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 12.44836933333333
$arr[11,13] = 37.345108
$arr[11,14] = 0.494893861707853
$arr[11,15] = 0.4948938617078529
$arr[11,16] = 747.287972310748
$arr[11,17] = 6725.591750796733
$arr[11,18] = 0.09420969391531678
$arr[11,19] = 0.09420969391531675

$arr[12,0] = "sCs"
$arr[12,1] = "Lpl"
$arr[12,2] = "Sdc1"
$arr[12,3] = "ECs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 16.00705666666667
$arr[12,7] = 48.02117
$arr[12,8] = 0.05075975073168155
$arr[12,9] = 0.05075975073168155
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 0.9159313333333333
$arr[12,13] = 2.747794
$arr[12,14] = 0.03641350786393945
$arr[12,15] = 0.03641350786393944
$arr[12,16] = 14.66136475544222
$arr[12,17] = 131.95228279898
$arr[12,18] = 0.001848340582439693
$arr[12,19] = 0.001848340582439692

$arr[13,0] = "sCs"
$arr[13,1] = "Lpl"
$arr[13,2] = "Sdc1"
$arr[13,3] = "FAPs"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 16.00705666666667
$arr[13,7] = 48.02117
$arr[13,8] = 0.05075975073168155
$arr[13,9] = 0.05075975073168155
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 10.56834633333333
$arr[13,13] = 31.705039
$arr[13,14] = 0.4201521973455823
$arr[13,15] = 0.4201521973455822
$arr[13,16] = 169.1681186306255
$arr[13,17] = 1522.51306767563
$arr[13,18] = 0.02132682080663003
$arr[13,19] = 0.02132682080663003

$arr[14,0] = "sCs"
$arr[14,1] = "Lpl"
$arr[14,2] = "Sdc1"
$arr[14,3] = "M2"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 16.00705666666667
$arr[14,7] = 48.02117
$arr[14,8] = 0.05075975073168155
$arr[14,9] = 0.05075975073168155
$arr[14,10] = 2
$arr[14,11] = 0.6666666666666666
$arr[14,12] = 1.220967333333333
$arr[14,13] = 3.662902
$arr[14,14] = 0.04854043308262539
$arr[14,15] = 0.04854043308262537
$arr[14,16] = 19.54409329281556
$arr[14,17] = 175.89683963534
$arr[14,18] = 0.002463900283681934
$arr[14,19] = 0.002463900283681932

$arr[15,0] = "sCs"
$arr[15,1] = "Lpl"
$arr[15,2] = "Sdc1"
$arr[15,3] = "sCs"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 16.00705666666667
$arr[15,7] = 48.02117
$arr[15,8] = 0.05075975073168155
$arr[15,9] = 0.05075975073168155
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 12.44836933333333
$arr[15,13] = 37.345108
$arr[15,14] = 0.494893861707853
$arr[15,15] = 0.4948938617078529
$arr[15,16] = 199.2617533262622
$arr[15,17] = 1793.35577993636
$arr[15,18] = 0.0251206890589299
$arr[15,19] = 0.02512068905892989

$ws.Range("A2:T17").Value = $arr
